$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header text (Volume number and date range) ---
$ws.Range("A8").Value = "Volume 29   Number  40"
$ws.Range("C9").Value = "Report Covering the Week  10/3/2022  Through  10/9/2022"

# --- Update data table cells (rows 15-30) ---
$ws.Range("C14").Copy($ws.Range("D15"))
$ws.Range("E14").Copy($ws.Range("E15"))
$ws.Range("C14").Copy($ws.Range("F15"))
$ws.Range("H15").Value = -100
$ws.Range("C16").Value = 8
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = 14.285714285714
$ws.Range("F16").Value = 16
$ws.Range("G16").Value = 20
$ws.Range("H16").Value = -20
$ws.Range("I16").Value = 122
$ws.Range("J16").Value = 103
$ws.Range("K16").Value = 18.446601941747
$ws.Range("L16").Value = 25.773195876288
$ws.Range("M16").Value = 52.5
$ws.Range("N16").Value = -74.845360824742
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 15
$ws.Range("G17").Value = 18
$ws.Range("H17").Value = -16.666666666666
$ws.Range("I17").Value = 103
$ws.Range("J17").Value = 100
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 27.160493827160
$ws.Range("M17").Value = 1.980198019801
$ws.Range("N17").Value = -46.907216494845
$ws.Range("I14").Copy($ws.Range("C18"))
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 16
$ws.Range("H18").Value = -62.5
$ws.Range("I18").Value = 135
$ws.Range("J18").Value = 96
$ws.Range("K18").Value = 40.625
$ws.Range("L18").Value = -0.735294117647
$ws.Range("M18").Value = 62.650602409638
$ws.Range("N18").Value = -71.875
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 17
$ws.Range("E19").Value = -23.529411764705
$ws.Range("F19").Value = 69
$ws.Range("G19").Value = 50
$ws.Range("H19").Value = 38
$ws.Range("I19").Value = 548
$ws.Range("J19").Value = 375
$ws.Range("K19").Value = 46.133333333333
$ws.Range("L19").Value = 31.100478468899
$ws.Range("M19").Value = 12.989690721649
$ws.Range("N19").Value = -14.774494556765
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = -57.142857142857
$ws.Range("I20").Value = 63
$ws.Range("J20").Value = 42
$ws.Range("K20").Value = 50
$ws.Range("L20").Value = 61.538461538461
$ws.Range("M20").Value = 80
$ws.Range("N20").Value = -83.887468030690
$ws.Range("C21").Value = 28
$ws.Range("D21").Value = 37
$ws.Range("E21").Value = -24.324324324324
$ws.Range("F21").Value = 109
$ws.Range("G21").Value = 113
$ws.Range("H21").Value = -3.539823008849
$ws.Range("I21").Value = 981
$ws.Range("J21").Value = 727
$ws.Range("K21").Value = 34.938101788170
$ws.Range("L21").Value = 25.769230769230
$ws.Range("M21").Value = 24.020227560050
$ws.Range("N21").Value = -55.790896800360
$ws.Range("C14").Copy($ws.Range("C22"))
$ws.Range("I14").Copy($ws.Range("D22"))
$ws.Range("D22").Value = 2
$ws.Range("K14").Copy($ws.Range("E22"))
$ws.Range("E22").Value = -100
$ws.Range("F22").Value = 5
$ws.Range("I14").Copy($ws.Range("G22"))
$ws.Range("G22").Value = 3
$ws.Range("K14").Copy($ws.Range("H22"))
$ws.Range("H22").Value = 66.666666666666
$ws.Range("I22").Value = 24
$ws.Range("J22").Value = 13
$ws.Range("K22").Value = 84.615384615384
$ws.Range("L22").Value = 118.181818181818
$ws.Range("M22").Value = 140
$ws.Range("G23").Value = 6
$ws.Range("H23").Value = -66.666666666666
$ws.Range("I23").Value = 44
$ws.Range("J23").Value = 52
$ws.Range("K23").Value = -15.384615384615
$ws.Range("L23").Value = 22.222222222222
$ws.Range("M23").Value = 29.411764705882
$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 14
$ws.Range("E24").Value = 35.714285714285
$ws.Range("F24").Value = 65
$ws.Range("G24").Value = 46
$ws.Range("H24").Value = 41.304347826087
$ws.Range("I24").Value = 582
$ws.Range("J24").Value = 413
$ws.Range("K24").Value = 40.920096852300
$ws.Range("L24").Value = 1.041666666666
$ws.Range("M24").Value = -13.521545319465
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = -12.5
$ws.Range("F25").Value = 34
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 298
$ws.Range("J25").Value = 246
$ws.Range("K25").Value = 21.138211382113
$ws.Range("L25").Value = 48.258706467661
$ws.Range("M25").Value = 12.452830188679
$ws.Range("C14").Copy($ws.Range("D26"))
$ws.Range("E14").Copy($ws.Range("E26"))
$ws.Range("C14").Copy($ws.Range("F26"))
$ws.Range("H26").Value = -100
$ws.Range("C14").Copy($ws.Range("C27"))
$ws.Range("I14").Copy($ws.Range("D27"))
$ws.Range("D27").Value = 4
$ws.Range("K14").Copy($ws.Range("E27"))
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 9
$ws.Range("H27").Value = -55.555555555555
$ws.Range("I27").Value = 46
$ws.Range("J27").Value = 40
$ws.Range("K27").Value = 15
$ws.Range("L27").Value = 48.387096774193
$ws.Range("C14").Copy($ws.Range("F30"))
$ws.Range("H30").Value = -100
$ws.Range("L30").Value = 15.384615384615
